$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 (ID 14)
$ws.Range("B17").Value = "Algunas de las actividades retrasadas no cuentan con seguimiento."
$ws.Range("C17").Value = "Ventas"
$ws.Range("D17").Value = 42482
$ws.Range("D17").NumberFormat = "m/d/yy"
$ws.Range("E17").Value = 42482
$ws.Range("E17").NumberFormat = "m/d/yy"
$ws.Range("F17").Value = "Cerrada"
$ws.Range("G17").Value = "Esta anomalia solo se reporta ya que son tareas ya realizadas. Se presentan para futuras actividades y se escalan."
$ws.Range("G17").WrapText = $true
$ws.Rows.Item(17).RowHeight = 75

# Row 18 (ID 15)
$ws.Range("B18").Value = "La auditoria a procesos y productos no ha sido completada."
$ws.Range("C18").Value = "Calidad"
$ws.Range("D18").Value = 42482
$ws.Range("D18").NumberFormat = "m/d/yy"
$ws.Range("E18").Value = 42482
$ws.Range("E18").NumberFormat = "m/d/yy"
$ws.Range("F18").Value = "Cerrada"
$ws.Range("G18").Value = "Esta anomalia solo se reporta ya que son tareas ya realizadas. Se presentan para futuras actividades y se escalan."
$ws.Range("G18").WrapText = $true
$ws.Rows.Item(18).RowHeight = 75

# Row 19 (ID 16)
$ws.Range("B19").Value = "La tarea de 20 y 21 no estan completadas."
$ws.Range("C19").Value = "Compras"
$ws.Range("D19").Value = 42482
$ws.Range("D19").NumberFormat = "m/d/yy"
$ws.Range("E19").Value = 42482
$ws.Range("E19").NumberFormat = "m/d/yy"
$ws.Range("F19").Value = "Cerrada"
$ws.Range("G19").Value = "Esta anomalia solo se reporta ya que son tareas ya realizadas. Se presentan para futuras actividades y se escalan."
$ws.Range("G19").WrapText = $true
$ws.Rows.Item(19).RowHeight = 75

# Row 20 (ID 17)
$ws.Range("B20").Value = "Faltan algunas tareas por completar."
$ws.Range("C20").Value = "Compras"
$ws.Range("D20").Value = 42482
$ws.Range("D20").NumberFormat = "m/d/yy"
$ws.Range("E20").Value = 42482
$ws.Range("E20").NumberFormat = "m/d/yy"
$ws.Range("F20").Value = "Cerrada"
$ws.Range("G20").Value = "Esta anomalia solo se reporta ya que son tareas ya realizadas. Se presentan para futuras actividades y se escalan."
$ws.Range("G20").WrapText = $true
$ws.Rows.Item(20).RowHeight = 75

# View state: scroll to show rows from 17, select D23
$ws.Activate()
$excel.ActiveWindow.TopLeftCell = $ws.Range("A17")
$ws.Range("D23").Select()

# Calculation properties: disable concurrent (multi-threaded) calculation
$excel.MultiThreadedCalculation.Enabled = $false
